$p = $ppt.ActivePresentation
Write-Output "Designs before:" $p.Designs.Count
$d2 = $p.Designs.Add()
Write-Output "Designs after add:" $p.Designs.Count
for ($i = 1; $i -le $p.Designs.Count; $i++) {
    $d = $p.Designs.Item($i)
    Write-Output "Design $i : $($d.Name)"
}
